# Trade #88 closed at 2026-02-16 21:38:01 - momentum DOWN +0.000%
#
# This script:
#  1. Updates the rollup metrics on "Summary" (rows 2-3) and "Comparison"
#     (row 2, leadlag) now that trade #57 (leadlag) has closed.
#  2. Marks trade #57 on the "leadlag" sheet (row 46) as CLOSED with its
#     exit price / P&L / exit reason.
#  3. Appends the newly opened trade #88 (momentum, DOWN) as a new row on
#     the "momentum" sheet.
#  4. Appends a copy of the now-closed trade #57 to the "All Trades" sheet
#     (it is rewritten there whenever a trade's state changes).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1a. Summary sheet - OVERALL (row 2) and leadlag STRATEGY (row 3) rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 57

$summary.Range("D2:F2").NumberFormat = "@"
$summary.Range("D2").Value = "66.7%"
$summary.Range("E2").Value = "+15.0972%"
$summary.Range("F2").Value = "+0.2649%"

$summary.Range("C3").Value = 66

$summary.Range("D3:F3").NumberFormat = "@"
$summary.Range("D3").Value = "42.4%"
$summary.Range("E3").Value = "+10.2464%"
$summary.Range("F3").Value = "+0.1552%"

# ---------------------------------------------------------------------
# 1b. Comparison sheet - leadlag row (row 2)
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 66

$comparison.Range("C2").NumberFormat = "@"
$comparison.Range("C2").Value = "42.4%"

$comparison.Range("D2").NumberFormat = "@"
$comparison.Range("D2").Value = "2.93"

$comparison.Range("F2").NumberFormat = "@"
$comparison.Range("F2").Value = "-0.3119%"

$comparison.Range("G2").NumberFormat = "@"
$comparison.Range("G2").Value = "1.78"

# ---------------------------------------------------------------------
# 2. leadlag sheet - close out trade #57 (row 46)
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G46").Value = 68846.535416
$leadlag.Range("H46").Value = "CLOSED"
$leadlag.Range("I46").Value = -0.0763
$leadlag.Range("J46").Value = -0.76
$leadlag.Range("M46").Value = "time_exit_5min"
$leadlag.Range("N46").Value = 5

# ---------------------------------------------------------------------
# 3. momentum sheet - append newly opened trade #88 (row 23)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("A23").Value = 88

$momentum.Range("B23:E23").NumberFormat = "@"
$momentum.Range("B23").Value = "2026-02-16"
$momentum.Range("C23").Value = "21:38:01"
$momentum.Range("D23").Value = "momentum"
$momentum.Range("E23").Value = "DOWN"

$momentum.Range("F23").Value = 68481.19500000001
# G23 stays an (empty) text cell - trade is still OPEN, no exit price yet.
$momentum.Range("G23").Value = "'"
$momentum.Range("H23").Value = "OPEN"
$momentum.Range("I23").Value = 0
$momentum.Range("J23").Value = 0
$momentum.Range("K23").Value = 0.9
$momentum.Range("L23").Value = "Downward momentum: -0.132% over 10 samples"
# M23 stays an (empty) text cell - no exit reason yet.
$momentum.Range("M23").Value = "'"
$momentum.Range("N23").Value = 0

# ---------------------------------------------------------------------
# 4. All Trades sheet - append the closed leadlag trade #57 (row 58)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A58").Value = 57

$allTrades.Range("B58:E58").NumberFormat = "@"
$allTrades.Range("B58").Value = "2026-02-16"
$allTrades.Range("C58").Value = "21:32:55"
$allTrades.Range("D58").Value = "leadlag"
$allTrades.Range("E58").Value = "DOWN"

$allTrades.Range("F58").Value = 68794.035
$allTrades.Range("G58").Value = 68846.535416
$allTrades.Range("H58").Value = "CLOSED"
$allTrades.Range("I58").Value = -0.07630000000000001
$allTrades.Range("J58").Value = -0.76
$allTrades.Range("K58").Value = 0.75
$allTrades.Range("L58").Value = "Binance leading with -0.114% move"
$allTrades.Range("M58").Value = "time_exit_5min"
$allTrades.Range("N58").Value = 5
